# personnel edits to fix roles, spacing, acknowledgments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Fix role spacing / correct role assignments
$ws.Cells.Item(5, 7).Value = "principal Investigator"   # Rachel Stanley
$ws.Cells.Item(7, 7).Value = "metadata Provider"         # Jaxine Wolfe
$ws.Cells.Item(8, 7).Value = "metadata Provider"         # Kate Morkeski
$ws.Cells.Item(9, 7).Value = "creator"                   # Arshia Mehta (was student)

# Make Personnel the active sheet, with G8 selected
$ws.Activate()
$ws.Range("G8").Select()
